$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1258.439
$ws.Range("I15").Value = 1258.439
$ws.Range("K15").Value = 3775.317
$ws.Range("M15").Value = -3606.317

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

# Row 132
$ws.Range("H132").Value = 2749194.5
$ws.Range("I132").Value = 2977935.5
$ws.Range("J132").Value = 4301.5
$ws.Range("K132").Value = 8933806.5
$ws.Range("L132").Value = 12904.5
$ws.Range("M132").Value = -8931276.5
$ws.Range("N132").Value = -17964.5

# Row 137
$ws.Range("H137").Value = 1453.7587
$ws.Range("I137").Value = 1398.2778
$ws.Range("J137").Value = 1544.5454
$ws.Range("K137").Value = 4194.8334
$ws.Range("L137").Value = 4633.6362
$ws.Range("M137").Value = -1644.8334
$ws.Range("N137").Value = -9733.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1159.4445
$ws.Range("I45").Value = 652.5
$ws.Range("J45").Value = 1565
$ws.Range("K45").Value = 652.5
$ws.Range("L45").Value = 1565
$ws.Range("M45").Value = -275.5
$ws.Range("N45").Value = -2319

# Row 101
$ws.Range("H101").Value = 27990.477
$ws.Range("J101").Value = 27990.477
$ws.Range("L101").Value = 27990.477
$ws.Range("N101").Value = -34480.477

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 342.75
$ws.Range("I80").Value = 1025
$ws.Range("J80").Value = 172.1875
$ws.Range("K80").Value = 1025
$ws.Range("L80").Value = 172.1875
$ws.Range("M80").Value = -27
$ws.Range("N80").Value = -2168.1875

# Row 83
$ws.Range("H83").Value = 342.75
$ws.Range("I83").Value = 1025
$ws.Range("J83").Value = 172.1875
$ws.Range("K83").Value = 5125
$ws.Range("L83").Value = 860.9375
$ws.Range("M83").Value = -133
$ws.Range("N83").Value = -10844.9375

# Row 86
$ws.Range("H86").Value = 2268.5
$ws.Range("I86").Value = 1437
$ws.Range("J86").Value = 3100
$ws.Range("K86").Value = 1437
$ws.Range("L86").Value = 3100
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -5346

# Row 89
$ws.Range("H89").Value = 2268.5
$ws.Range("I89").Value = 1437
$ws.Range("J89").Value = 3100
$ws.Range("K89").Value = 7185
$ws.Range("L89").Value = 15500
$ws.Range("M89").Value = -1569
$ws.Range("N89").Value = -26732

# Row 107
$ws.Range("H107").Value = 1190.6086
$ws.Range("I107").Value = 707.2
$ws.Range("J107").Value = 1562.4615
$ws.Range("K107").Value = 707.2
$ws.Range("L107").Value = 1562.4615
$ws.Range("M107").Value = 1212.8
$ws.Range("N107").Value = -5402.4615

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5129422.5
$ws.Range("I31").Value = 1293.2307
$ws.Range("J31").Value = 12821617
$ws.Range("K31").Value = 1293.2307
$ws.Range("L31").Value = 12821617
$ws.Range("M31").Value = -998.2307
$ws.Range("N31").Value = -12822207

# Row 34
$ws.Range("H34").Value = 5129422.5
$ws.Range("I34").Value = 1293.2307
$ws.Range("J34").Value = 12821617
$ws.Range("K34").Value = 1293.2307
$ws.Range("L34").Value = 12821617
$ws.Range("M34").Value = -1091.2307
$ws.Range("N34").Value = -12822021

# Row 58
$ws.Range("H58").Value = 960.3125
$ws.Range("I58").Value = 823.3333
$ws.Range("J58").Value = 1700
$ws.Range("K58").Value = 823.3333
$ws.Range("L58").Value = 1700
$ws.Range("M58").Value = -620.3333
$ws.Range("N58").Value = -2106

# Row 99
$ws.Range("H99").Value = 2639.25
$ws.Range("I99").Value = 2100
$ws.Range("J99").Value = 2716.2856
$ws.Range("K99").Value = 2100
$ws.Range("L99").Value = 2716.2856
$ws.Range("M99").Value = -602
$ws.Range("N99").Value = -5712.2856

# Row 109
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080

# Row 126
$ws.Range("H126").Value = 2639.25
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 2716.2856
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 8148.8568
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -13088.8568

# Row 134
$ws.Range("H134").Value = 699.1111
$ws.Range("I134").Value = 503.1
$ws.Range("J134").Value = 1259.1428
$ws.Range("K134").Value = 1509.3
$ws.Range("L134").Value = 3777.4284
$ws.Range("M134").Value = 1025.7
$ws.Range("N134").Value = -8847.4284

# Row 136
$ws.Range("H136").Value = 960.3125
$ws.Range("I136").Value = 823.3333
$ws.Range("J136").Value = 1700
$ws.Range("K136").Value = 2469.9999
$ws.Range("L136").Value = 5100
$ws.Range("M136").Value = 80.0001
$ws.Range("N136").Value = -10200

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 113.23077
$ws.Range("I23").Value = 21.5
$ws.Range("J23").Value = 129.90909
$ws.Range("K23").Value = 64.5
$ws.Range("L23").Value = 389.72727
$ws.Range("M23").Value = 170.5
$ws.Range("N23").Value = -859.72727

# Row 38
$ws.Range("H38").Value = 18.714285
$ws.Range("I38").Value = 8.2
$ws.Range("J38").Value = 24.555555
$ws.Range("K38").Value = 24.6
$ws.Range("L38").Value = 73.666665
$ws.Range("M38").Value = 322.4
$ws.Range("N38").Value = -767.666665

# Row 62
$ws.Range("H62").Value = 3150
$ws.Range("J62").Value = 3150
$ws.Range("L62").Value = 9450
$ws.Range("N62").Value = -10822

# Row 63
$ws.Range("H63").Value = 11731.393
$ws.Range("I63").Value = 3622.375
$ws.Range("J63").Value = 14975
$ws.Range("K63").Value = 10867.125
$ws.Range("L63").Value = 44925
$ws.Range("M63").Value = -10118.125
$ws.Range("N63").Value = -46423

# Row 64
$ws.Range("H64").Value = 880
$ws.Range("I64").Value = 880
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 2640
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2370
$ws.Range("N64").ClearContents()

# Row 65
$ws.Range("H65").Value = 3150
$ws.Range("J65").Value = 3150
$ws.Range("L65").Value = 28350
$ws.Range("N65").Value = -35214

# Row 66
$ws.Range("H66").Value = 11731.393
$ws.Range("I66").Value = 3622.375
$ws.Range("J66").Value = 14975
$ws.Range("K66").Value = 32601.375
$ws.Range("L66").Value = 134775
$ws.Range("M66").Value = -28857.375
$ws.Range("N66").Value = -142263

# Row 67
$ws.Range("H67").Value = 880
$ws.Range("I67").Value = 880
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 2640
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -1704
$ws.Range("N67").ClearContents()

# Row 107
$ws.Range("H107").Value = 500.2143
$ws.Range("I107").Value = 357.9
$ws.Range("J107").Value = 579.2778
$ws.Range("K107").Value = 1073.7
$ws.Range("L107").Value = 1737.8334
$ws.Range("M107").Value = 846.3
$ws.Range("N107").Value = -5577.8334

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2548.7
$ws.Range("I7").Value = 2124.8
$ws.Range("J7").Value = 2690
$ws.Range("K7").Value = 2124.8
$ws.Range("L7").Value = 2690
$ws.Range("M7").Value = -2012.8
$ws.Range("N7").Value = -2914

# Row 126
$ws.Range("H126").Value = 2548.7
$ws.Range("I126").Value = 2124.8
$ws.Range("J126").Value = 2690
$ws.Range("K126").Value = 6374.4
$ws.Range("L126").Value = 8070
$ws.Range("M126").Value = -3904.4
$ws.Range("N126").Value = -13010

# Row 132
$ws.Range("H132").Value = 5101.735
$ws.Range("I132").Value = 7644.778
$ws.Range("J132").Value = 1980.7273
$ws.Range("K132").Value = 22934.334
$ws.Range("L132").Value = 5942.1819
$ws.Range("M132").Value = -20404.334
$ws.Range("N132").Value = -11002.1819

$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 20000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 20000
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -22774

Write-Output "Applied all Garuda_Profits updates"